# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 ("SOURCES OF FINANCE") gets a new gallery table
#    style applied (Table Tools > Design > Table Styles gallery):
#      {26E61C38-7681-4EC0-AA67-B9C4772F33E8}  ->  {E6A518CE-1F66-44AD-A7C1-02022D944329}
#
# 2) The deck's theme colour palette is changed from the custom "Integral"
#    palette back to the stock Office palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), e.g. via Design > Variants > Colors > "Office".
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 ---------------------------------------

$targetStyleId = "{E6A518CE-1F66-44AD-A7C1-02022D944329}"

$tableShape = $null
$slideWithTable = $null
for ($si = 1; $si -le $p.Slides.Count -and $tableShape -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $candidate = $slide.Shapes.Item($shi)
        if ($candidate.HasTable) {
            $tableShape = $candidate
            $slideWithTable = $slide
            break
        }
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle($targetStyleId)
}

# --- 2. Swap the theme colour scheme back to the stock "Office" palette -----
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink — RGB() uses PowerPoint's
# BGR-packed long, so 0xBBGGRR, not 0xRRGGBB.)

$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

$firstSlide = $p.Slides.Item(1)
$themeColors = $firstSlide.ThemeColorScheme
for ($ci = 1; $ci -le $officeColors.Length; $ci++) {
    $themeColors.Colors($ci).RGB = $officeColors[$ci - 1]
}
